# Apply the "Fit QRF models updated script" edit to ModelCovSelected.xlsx
# Sheets (1-based Worksheets index):
#   4 -> CHaMP_Redds_Steelhead
#   5 -> CHaMP_Winter_Chinook
#   6 -> CHaMP_Winter_Steelhead

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet "CHaMP_Redds_Steelhead": reset scroll position (remove topLeftCell="A10")
# -------------------------------------------------------------------
$wsReddsSteelhead = $wb.Worksheets.Item(4)
$wsReddsSteelhead.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# -------------------------------------------------------------------
# Sheet "CHaMP_Winter_Steelhead": data updates
# -------------------------------------------------------------------
$wsWinterSteelhead = $wb.Worksheets.Item(6)
$wsWinterSteelhead.Activate()

# Row 9 (Ucut_Length): QRF2_trimmed flag 1 -> 0
$wsWinterSteelhead.Cells.Item(9, 4).Value = 0

# Row 13: metric changed from Discharge_fish to Q, QRF2_trimmed flag 1 -> 0
$wsWinterSteelhead.Cells.Item(13, 2).Value = "Q"
$wsWinterSteelhead.Cells.Item(13, 4).Value = 0

# Row 19 (SubEstCbl): QRF2_trimmed flag 0 -> 1
$wsWinterSteelhead.Cells.Item(19, 4).Value = 1

# Row 21 (SubEstBldr): QRF2_trimmed flag 0 -> 1
$wsWinterSteelhead.Cells.Item(21, 4).Value = 1

# Row 22 (SubEstCandBldr): QRF2_trimmed flag 1 -> 0, note text updated
$wsWinterSteelhead.Cells.Item(22, 4).Value = 0
$wsWinterSteelhead.Cells.Item(22, 5).Value = "Remove to be consistant with other models"

# Column widths: C and D get wider / custom widths
$wsWinterSteelhead.Columns.Item(3).ColumnWidth = 12.28
$wsWinterSteelhead.Columns.Item(4).ColumnWidth = 26.83

# Selection moves to D26, sheet no longer the selected tab (handled by later activations below)
$wsWinterSteelhead.Range("D26").Select()

# -------------------------------------------------------------------
# Sheet "CHaMP_Winter_Chinook": data updates (becomes the active tab)
# -------------------------------------------------------------------
$wsWinterChinook = $wb.Worksheets.Item(5)
$wsWinterChinook.Activate()

# Row 14: metric changed from Discharge_fish to Q, QRF2_trimmed flag 1 -> 0, note removed
$wsWinterChinook.Cells.Item(14, 2).Value = "Q"
$wsWinterChinook.Cells.Item(14, 4).Value = 0
$wsWinterChinook.Cells.Item(14, 5).Clear()

# Row 20 (SubEstBldr): QRF2_trimmed flag 0 -> 1
$wsWinterChinook.Cells.Item(20, 4).Value = 1

# Row 21 (SubEstCbl): QRF2_trimmed flag 0 -> 1
$wsWinterChinook.Cells.Item(21, 4).Value = 1

# Row 22 (SubEstCandBldr): QRF2_trimmed flag 1 -> 0, note added
$wsWinterChinook.Cells.Item(22, 4).Value = 0
$wsWinterChinook.Cells.Item(22, 5).Value = "Remove to be consistant with other models"

# Selection moves to D33, and this sheet becomes the active/selected tab
$wsWinterChinook.Range("D33").Select()
$wsWinterChinook.Activate()
